$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the Price (D) and Volume(1h) (E) columns so that
# numeric-looking strings (e.g. "0.999", "3.24", "0.0414") are preserved verbatim
# as text instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# U+2083 SUBSCRIPT THREE, used in the PEPE price text "0.0<sub3>0788".
$subThree = [string][char]0x2083

$ws.Range("D2").Value = '65.565.05'
$ws.Range("E2").Value = '  -3.02%  '
$ws.Range("D3").Value = '3.496.30'
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '554.60'
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("D6").Value = '178.86'
$ws.Range("E6").Value = '  -5.18%  '
$ws.Range("E7").Value = '  +3.83%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '0.632'
$ws.Range("E9").Value = '  -1.10%  '
$ws.Range("E10").Value = '  +1.94%  '
$ws.Range("D11").Value = '53.78'
$ws.Range("E11").Value = '  -5.09%  '
$ws.Range("E12").Value = '  -1.55%  '
$ws.Range("D13").Value = '9.23'
$ws.Range("E13").Value = '  -2.51%  '
$ws.Range("D14").Value = '4.052.58'
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("D15").Value = '3.490.39'
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '18.46'
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = '0.121'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").Value = '12.07'
$ws.Range("E18").Value = '  +1.96%  '
$ws.Range("D19").Value = '65.507.58'
$ws.Range("E19").Value = '  -3.46%  '
$ws.Range("D20").Value = '0.991'
$ws.Range("E20").Value = '  -1.73%  '
$ws.Range("D21").Value = '417.80'
$ws.Range("E21").Value = '  +3.25%  '
$ws.Range("D22").Value = '4.04'
$ws.Range("E22").Value = '  +1.99%  '
$ws.Range("D23").Value = '86.46'
$ws.Range("E23").Value = '  +2.13%  '
$ws.Range("D24").Value = '4.12'
$ws.Range("E24").Value = '  -2.11%  '
$ws.Range("D25").Value = '12.77'
$ws.Range("E25").Value = '  +7.69%  '
$ws.Range("D26").Value = '10.79'
$ws.Range("E26").Value = '  -10.98%  '
$ws.Range("D27").Value = '2.85'
$ws.Range("E27").Value = '  -2.76%  '
$ws.Range("E28").Value = '  -3.28%  '
$ws.Range("D29").Value = '9.04'
$ws.Range("E29").Value = '  +4.56%  '
$ws.Range("D30").Value = '30.27'
$ws.Range("E30").Value = '  -0.26%  '
$ws.Range("D31").Value = '6.50'
$ws.Range("D32").Value = '608.48'
$ws.Range("E32").Value = '  -11.26%  '
$ws.Range("D33").Value = '11.72'
$ws.Range("E33").Value = '  +0.13%  '
$ws.Range("D34").Value = '0.110'
$ws.Range("E34").Value = '  -0.79%  '
$ws.Range("D35").Value = '59.55'
$ws.Range("E35").Value = '  -0.50%  '
$ws.Range("E36").Value = '  +9.72%  '
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").Value = '37.40'
$ws.Range("E38").Value = '  -3.98%  '
$ws.Range("D39").Value = '0.0' + $subThree + '0788'
$ws.Range("E39").Value = '  -5.97%  '
$ws.Range("D40").Value = '3.378.90'
$ws.Range("E40").Value = '  +10.84%  '
$ws.Range("D41").Value = '0.381'
$ws.Range("E41").Value = '  -5.77%  '
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").Value = '3.25'
$ws.Range("E43").Value = '  -5.73%  '
$ws.Range("D44").Value = '2.85'
$ws.Range("E44").Value = '  -3.08%  '
$ws.Range("D45").Value = '2.54'
$ws.Range("E45").Value = '  -9.42%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '0.0414'
$ws.Range("E46").Value = '  -1.97%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '3.24'
$ws.Range("E47").Value = '  -0.97%  '
$ws.Range("D48").Value = '2.70'
$ws.Range("E48").Value = '  -1.78%  '
$ws.Range("E49").Value = '  +1.25%  '
$ws.Range("D50").Value = '8.49'
$ws.Range("E50").Value = '  -3.82%  '
$ws.Range("D51").Value = '137.86'
$ws.Range("E51").Value = '  -0.85%  '
